# Auto-generated Excel COM-interop edit script
# Reproduces the numeric (F/G column) corrections plus the
# newly-added "生如夏花国乙only" event row inserted into the
# "全部类型" aggregate sheet (rows 38-45 content shift).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1111
$ws.Range("F4").Value = 191
$ws.Range("F6").Value = 1747
$ws.Range("F7").Value = 645
$ws.Range("F9").Value = 402
$ws.Range("F10").Value = 4060
$ws.Range("F11").Value = 41
$ws.Range("F14").Value = 971
$ws.Range("F18").Value = 2923
$ws.Range("F19").Value = 1749
$ws.Range("F24").Value = 901
$ws.Range("F26").Value = 2172
$ws.Range("F28").Value = 2259
$ws.Range("F30").Value = 656
$ws.Range("F31").Value = 426
$ws.Range("F34").Value = 385
$ws.Range("F35").Value = 1033
$ws.Range("F36").Value = 854
$ws.Range("F38").Value = 285
$ws.Range("F39").Value = 488
$ws.Range("F40").Value = 340

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 12
$ws.Range("F10").Value = 864
$ws.Range("F13").Value = 2
$ws.Range("F15").Value = 13

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1111
$ws.Range("F4").Value = 191
$ws.Range("F7").Value = 1747
$ws.Range("F8").Value = 645
$ws.Range("F10").Value = 402
$ws.Range("F11").Value = 4060
$ws.Range("F12").Value = 41
$ws.Range("F13").Value = 12
$ws.Range("F19").Value = 2923
$ws.Range("F21").Value = 1749
$ws.Range("F26").Value = 864
$ws.Range("F28").Value = 2
$ws.Range("F29").Value = 901
$ws.Range("F30").Value = 2172
$ws.Range("F31").Value = 13
$ws.Range("F34").Value = 2259
$ws.Range("F35").Value = 656
$ws.Range("F36").Value = 426
$ws.Range("C38").Value = "杭州·生如夏花国乙only·日夜场"
$ws.Range("D38").Value = "祥符街道花园岗街181号 格乐利雅婚礼艺术中心(天空之城店)"
$ws.Range("E38").Value = "2024.07.20 10:00-07.20 22:30"
$ws.Range("F38").Value = 385
$ws.Range("G38").Value = 105
$ws.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=85496"
$ws.Range("I38").Value = "//i1.hdslb.com/bfs/openplatform/202405/Qut2ZdAi1715411977772.jpeg"
$ws.Range("B39").Value = "'2024-07-20"
$ws.Range("C39").Value = "杭州·第五届华盟次元嘉年华&周年庆狂欢"
$ws.Range("D39").Value = "创意路1号 中国智谷富春园区"
$ws.Range("E39").Value = "2024.07.20 10:00-07.21 17:00"
$ws.Range("F39").Value = 1033
$ws.Range("G39").Value = 58
$ws.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=84762"
$ws.Range("I39").Value = "//i0.hdslb.com/bfs/openplatform/202404/uE6OVg6T1713885553204.jpeg"
$ws.Range("B40").Value = "'2024-07-27"
$ws.Range("C40").Value = "杭州·夏之誓国乙only-日夜场"
$ws.Range("D40").Value = "北干街道萧杭路689号 杭州时尚外滩艺术中心"
$ws.Range("E40").Value = "2024.07.27 10:00-07.27 21:00"
$ws.Range("F40").Value = 854
$ws.Range("G40").Value = 69
$ws.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=83589"
$ws.Range("I40").Value = "//i2.hdslb.com/bfs/openplatform/202405/99kWb2dy1714964533903.png"
$ws.Range("B41").Value = "'2024-07-28"
$ws.Range("C41").Value = "杭州·火影忍者only"
$ws.Range("D41").Value = "康候圣街99号 顺丰创新中心"
$ws.Range("E41").Value = "2024.07.28 09:00-07.28 18:00"
$ws.Range("F41").Value = 1127
$ws.Range("G41").Value = 75
$ws.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=84243"
$ws.Range("I41").Value = "//i1.hdslb.com/bfs/openplatform/202405/isG309e51715657222196.jpeg"
$ws.Range("C42").Value = "杭州·AP动漫游戏嘉年华"
$ws.Range("D42").Value = "沈半路171号 Tcar汽车文化主题公园"
$ws.Range("E42").Value = "2024.08.03 09:00-08.04 17:00"
$ws.Range("F42").Value = 285
$ws.Range("G42").Value = 70
$ws.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=85527"
$ws.Range("I42").Value = "//i2.hdslb.com/bfs/openplatform/202405/JbVl16OE1715676665714.jpeg"
$ws.Range("B43").Value = "'2024-08-03"
$ws.Range("C43").Value = "杭州·梦漫星河动漫展"
$ws.Range("D43").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$ws.Range("E43").Value = "2024.08.03 10:00-08.04 17:00"
$ws.Range("F43").Value = 488
$ws.Range("G43").Value = 68
$ws.Range("H43").Value = "https://show.bilibili.com/platform/detail.html?id=82836"
$ws.Range("I43").Value = "//i0.hdslb.com/bfs/openplatform/202403/VFfQUJdD1711700169290.jpeg"
$ws.Range("C44").Value = "杭州·【七夕巨献·早鸟6折】真的爱你”致敬Beyond·黄家驹31周年演唱会·630乐团再现91殿堂级演出"
$ws.Range("D44").Value = "湖墅南路136-138号 浙话艺术剧院"
$ws.Range("E44").Value = "2024.08.10 19:30-08.10 21:30"
$ws.Range("F44").Value = 0
$ws.Range("H44").Value = "https://show.bilibili.com/platform/detail.html?id=85333"
$ws.Range("I44").Value = "//i1.hdslb.com/bfs/openplatform/202405/uYt32zt21715221330023.jpeg"
$ws.Range("C45").Value = "杭州·原神X星铁X绝区零only"
$ws.Range("D45").Value = "望江东路333号 杭州瑞莱克斯大酒店"
$ws.Range("E45").Value = "2024.08.10 10:00-08.10 17:00"
$ws.Range("F45").Value = 340
$ws.Range("G45").Value = 60
$ws.Range("H45").Value = "https://show.bilibili.com/platform/detail.html?id=82754"
$ws.Range("I45").Value = "//i1.hdslb.com/bfs/openplatform/202403/qA0LNJuF1710234461030.jpeg"
